$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.722.51"
$ws.Range("E2").Value = "  +2.33%  "
$ws.Range("D3").Value = "2.158.77"
$ws.Range("E3").Value = "  +2.62%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'227.07"
$ws.Range("E5").Value = "  -0.22%  "
$ws.Range("D6").Value = "'0.627"
$ws.Range("E6").Value = "  +1.57%  "
$ws.Range("D7").Value = "'63.27"
$ws.Range("E7").Value = "  +1.93%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +0.68%  "
$ws.Range("E10").Value = "  +0.46%  "
$ws.Range("E11").Value = "  -0.49%  "
$ws.Range("D12").Value = "'15.89"
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("D13").Value = "2.478.24"
$ws.Range("E13").Value = "  +2.48%  "
$ws.Range("D14").Value = "'21.84"
$ws.Range("E14").Value = "  -0.62%  "
$ws.Range("D15").Value = "'0.803"
$ws.Range("E15").Value = "  -0.31%  "
$ws.Range("D16").Value = "'5.47"
$ws.Range("E16").Value = "  -0.79%  "
$ws.Range("D17").Value = "2.156.46"
$ws.Range("E17").Value = "  +2.72%  "
$ws.Range("D18").Value = "39.597.84"
$ws.Range("E18").Value = "  +2.05%  "
$ws.Range("D19").Value = "'71.63"
$ws.Range("E19").Value = "  -0.11%  "
$ws.Range("E20").Value = "  -0.56%  "
$ws.Range("E21").Value = "  -0.59%  "
$ws.Range("D22").Value = "'230.03"
$ws.Range("E22").Value = "  +0.93%  "
$ws.Range("D24").Value = "'2.37"
$ws.Range("E24").Value = "  +2.38%  "
$ws.Range("D25").Value = "'2.31"
$ws.Range("E25").Value = "  -1.44%  "
$ws.Range("D26").Value = "'172.27"
$ws.Range("E26").Value = "  +0.35%  "
$ws.Range("D27").Value = "'9.54"
$ws.Range("E27").Value = "  -1.81%  "
$ws.Range("E28").Value = "  +2.39%  "
$ws.Range("D29").Value = "'1.45"
$ws.Range("E29").Value = "  +2.28%  "
$ws.Range("E30").Value = "  +2.58%  "
$ws.Range("D31").Value = "'2.69"
$ws.Range("E31").Value = "  +5.45%  "
$ws.Range("E32").Value = "  +1.36%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").Value = "'4.67"
$ws.Range("E34").Value = "  -1.16%  "
$ws.Range("D35").Value = "'6.94"
$ws.Range("E35").Value = "  -3.99%  "
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("E37").Value = "  +0.33%  "
$ws.Range("E38").Value = "  +2.80%  "
$ws.Range("D39").Value = "'5.11"
$ws.Range("E39").Value = "  +22.98%  "
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("D41").Value = "'102.72"
$ws.Range("E41").Value = "  +0.97%  "
$ws.Range("E42").Value = "  -0.72%  "
$ws.Range("D43").Value = "'17.60"
$ws.Range("E43").Value = "  -2.33%  "
$ws.Range("D44").Value = "1.515.59"
$ws.Range("E44").Value = "  -0.62%  "
$ws.Range("E45").Value = "  +0.57%  "
$ws.Range("D46").Value = "'2.80"
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("D47").Value = "'0.0920"
$ws.Range("E47").Value = "  +0.49%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "'1.09"
$ws.Range("E48").Value = "  +0.06%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "'7.75"
$ws.Range("E49").Value = "  -0.70%  "
$ws.Range("D50").Value = "'50.20"
$ws.Range("E50").Value = "  +9.05%  "
$ws.Range("E51").Value = "  +0.79%  "
